$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = '71.715.27'
$ws.Range("E2").Value = '  +3.46%  '
$ws.Range("D3").Value = '3.701.95'
$ws.Range("E3").Value = '  +8.27%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").Value = "'590.74"
$ws.Range("E5").Value = '  +1.54%  '
$ws.Range("D6").Value = "'180.71"
$ws.Range("E6").Value = '  +1.44%  '
$ws.Range("D7").Value = '3.683.78'
$ws.Range("E7").Value = '  +7.95%  '
$ws.Range("E8").Value = '  +4.28%  '
$ws.Range("E9").Value = '  +0.10%  '
$ws.Range("E10").Value = '  +2.78%  '
$ws.Range("D11").Value = "'0.614"
$ws.Range("E11").Value = '  +5.06%  '
$ws.Range("D12").Value = "'50.00"
$ws.Range("E12").Value = '  +3.72%  '
$ws.Range("E13").Value = '  +2.80%  '
$ws.Range("D14").Value = '4.288.11'
$ws.Range("E14").Value = '  +8.05%  '
$ws.Range("D15").Value = "'685.50"
$ws.Range("E15").Value = '  +0.91%  '
$ws.Range("D16").Value = "'9.05"
$ws.Range("E16").Value = '  +4.76%  '
$ws.Range("D17").Value = '3.707.76'
$ws.Range("E17").Value = '  +8.55%  '
$ws.Range("D18").Value = '71.775.16'
$ws.Range("E18").Value = '  +3.48%  '
$ws.Range("E19").Value = '  +2.26%  '
$ws.Range("D20").Value = "'18.18"
$ws.Range("E20").Value = '  +2.35%  '
$ws.Range("D21").Value = "'11.70"
$ws.Range("E21").Value = '  +3.40%  '
$ws.Range("D22").Value = "'6.49"
$ws.Range("E22").Value = '  +20.94%  '
$ws.Range("D23").Value = "'0.946"
$ws.Range("E23").Value = '  +3.97%  '
$ws.Range("D24").Value = "'17.85"
$ws.Range("E24").Value = '  +4.93%  '
$ws.Range("D25").Value = "'103.96"
$ws.Range("E25").Value = '  +3.13%  '
$ws.Range("D26").Value = "'4.06"
$ws.Range("E26").Value = '  +4.07%  '
$ws.Range("D27").Value = "'2.84"
$ws.Range("E27").Value = '  +5.58%  '
$ws.Range("D28").Value = "'10.27"
$ws.Range("E28").Value = '  +6.60%  '
$ws.Range("D29").Value = "'35.59"
$ws.Range("E29").Value = '  +6.17%  '
$ws.Range("E30").Value = '  +6.07%  '
$ws.Range("D31").Value = "'7.38"
$ws.Range("E31").Value = '  +7.48%  '
$ws.Range("D32").Value = "'4.27"
$ws.Range("E32").Value = '  +15.22%  '
$ws.Range("D33").Value = "'11.32"
$ws.Range("E33").Value = '  +2.80%  '
$ws.Range("D34").Value = "'564.95"
$ws.Range("E34").Value = '  +1.17%  '
$ws.Range("E35").Value = '  +4.24%  '
$ws.Range("D36").Value = "'59.73"
$ws.Range("E36").Value = '  +2.85%  '
$ws.Range("D37").Value = '3.755.79'
$ws.Range("E37").Value = '  +4.08%  '
$ws.Range("E38").Value = '  -0.14%  '
$ws.Range("D39").Value = "'0.145"
$ws.Range("E39").Value = '  +3.09%  '
$ws.Range("D40").Value = '0.0₃0781'
$ws.Range("E40").Value = '  +6.32%  '
$ws.Range("D41").Value = "'35.72"
$ws.Range("E41").Value = '  +2.32%  '
$ws.Range("E42").Value = '  +6.27%  '
$ws.Range("D43").Value = "'0.0467"
$ws.Range("E43").Value = '  +10.28%  '
$ws.Range("D44").Value = "'2.81"
$ws.Range("E44").Value = '  +4.55%  '
$ws.Range("D45").Value = "'0.353"
$ws.Range("E45").Value = '  +5.77%  '
$ws.Range("D46").Value = "'2.89"
$ws.Range("E46").Value = '  +8.67%  '
$ws.Range("E47").Value = '  -0.49%  '
$ws.Range("E48").Value = '  +4.03%  '
$ws.Range("E49").Value = '  +3.22%  '
$ws.Range("D50").Value = "'0.999"
$ws.Range("E50").Value = '  -0.05%  '
$ws.Range("D51").Value = "'135.66"
$ws.Range("E51").Value = '  +3.55%  '
